$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '2026-02-17 04:48:33'
$ws.Range("H2").Value = '44%'
$ws.Range("E3").Value = '2026-02-17 04:48:35'
$ws.Range("O3").Value = '-7.2 °C'
$ws.Range("E4").Value = '2026-02-17 04:48:37'
$ws.Range("J4").Value = '1015.2 hPa'
$ws.Range("O4").Value = '5.9 °C'
$ws.Range("E5").Value = '2026-02-17 04:48:40'
$ws.Range("E6").Value = '2026-02-17 04:48:42'
$ws.Range("H6").Value = '86%'
$ws.Range("J6").Value = '1014.8 hPa'
$ws.Range("E7").Value = '2026-02-17 04:48:45'
$ws.Range("J7").Value = '1014.4 hPa'
$ws.Range("N7").Value = '12.9 °C 4:29 TU'
$ws.Range("O7").Value = '14.4 °C'
$ws.Range("E8").Value = '2026-02-17 04:48:47'
$ws.Range("J8").Value = '1014.6 hPa'
$ws.Range("N8").Value = '8.9 °C 4:25 TU'
$ws.Range("O8").Value = '10.1 °C'
$ws.Range("E9").Value = '2026-02-17 04:48:49'
$ws.Range("O9").Value = '12.5 °C'
$ws.Range("E10").Value = '2026-02-17 04:48:52'
$ws.Range("H10").Value = '91%'
$ws.Range("L10").Value = '11.5 km/h - 32º 4:15 TU'
$ws.Range("M10").Value = '10.9 °C 4:19 TU'
$ws.Range("O10").Value = '7.8 °C'
$ws.Range("E11").Value = '2026-02-17 04:48:54'
$ws.Range("H11").Value = '35%'
$ws.Range("N11").Value = '2.8 °C 4:27 TU'
$ws.Range("O11").Value = '7.1 °C'
$ws.Range("E12").Value = '2026-02-17 04:48:56'
$ws.Range("O12").Value = '13.0 °C'
$ws.Range("E13").Value = '2026-02-17 04:48:58'
$ws.Range("H13").Value = '44%'
$ws.Range("J13").Value = '1016.4 hPa'
$ws.Range("N13").Value = '1.0 °C 4:12 TU'
$ws.Range("O13").Value = '4.3 °C'
$ws.Range("E14").Value = '2026-02-17 04:49:01'
$ws.Range("H14").Value = '59%'
$ws.Range("N14").Value = '12.4 °C 4:28 TU'
$ws.Range("O14").Value = '13.5 °C'
$ws.Range("E15").Value = '2026-02-17 04:49:03'
$ws.Range("N15").Value = '11.5 °C 4:02 TU'
$ws.Range("E16").Value = '2026-02-17 04:49:05'
$ws.Range("M16").Value = '-4.8 °C 4:05 TU'
$ws.Range("E17").Value = '2026-02-17 04:49:08'
$ws.Range("H17").Value = '52%'
$ws.Range("E18").Value = '2026-02-17 04:49:10'
$ws.Range("J18").Value = '1015.1 hPa'
$ws.Range("E19").Value = '2026-02-17 04:49:12'
$ws.Range("H19").Value = '68%'
$ws.Range("L19").Value = '31.0 km/h - 315º 4:23 TU'
$ws.Range("O19").Value = '6.2 °C'
$ws.Range("E20").Value = '2026-02-17 04:49:14'
$ws.Range("H20").Value = '41%'
$ws.Range("E21").Value = '2026-02-17 04:49:17'
$ws.Range("H21").Value = '30%'
$ws.Range("N21").Value = '7.3 °C 4:18 TU'
$ws.Range("E22").Value = '2026-02-17 04:49:19'
$ws.Range("E23").Value = '2026-02-17 04:49:21'
$ws.Range("H23").Value = '60%'
$ws.Range("E24").Value = '2026-02-17 04:49:24'
$ws.Range("N24").Value = '9.5 °C 4:25 TU'
$ws.Range("O24").Value = '9.9 °C'
$ws.Range("E25").Value = '2026-02-17 04:49:26'
$ws.Range("N25").Value = '-4.5 °C 4:11 TU'
$ws.Range("O25").Value = '-3.5 °C'
$ws.Range("E26").Value = '2026-02-17 04:49:28'
$ws.Range("E27").Value = '2026-02-17 04:49:31'
$ws.Range("H27").Value = '43%'
$ws.Range("N27").Value = '-3.2 °C 4:05 TU'
$ws.Range("E28").Value = '2026-02-17 04:49:33'
$ws.Range("H28").Value = '91%'
$ws.Range("J28").Value = '1015.3 hPa'
$ws.Range("E29").Value = '2026-02-17 04:49:35'
$ws.Range("E30").Value = '2026-02-17 04:49:38'
$ws.Range("J30").Value = '1014.2 hPa'
$ws.Range("N30").Value = '10.5 °C 4:29 TU'
$ws.Range("O30").Value = '12.2 °C'
$ws.Range("E31").Value = '2026-02-17 04:49:40'
$ws.Range("H31").Value = '62%'
$ws.Range("I31").Value = '0.1 mm'
$ws.Range("J31").Value = '1014.9 hPa'
$ws.Range("N31").Value = '7.7 °C 4:29 TU'
$ws.Range("O31").Value = '10.1 °C'
$ws.Range("E32").Value = '2026-02-17 04:49:42'
$ws.Range("H32").Value = '74%'
$ws.Range("N32").Value = '5.5 °C 4:21 TU'
$ws.Range("O32").Value = '6.5 °C'
$ws.Range("E33").Value = '2026-02-17 04:49:44'
$ws.Range("H33").Value = '40%'
$ws.Range("J33").Value = '1015.6 hPa'
$ws.Range("E34").Value = '2026-02-17 04:49:47'
$ws.Range("L34").Value = '70.2 km/h - 48º 4:06 TU'
$ws.Range("E35").Value = '2026-02-17 04:49:49'
$ws.Range("I35").Value = '2.3 mm'
$ws.Range("J35").Value = '1018.0 hPa'
$ws.Range("O35").Value = '5.5 °C'
$ws.Range("E36").Value = '2026-02-17 04:49:52'
$ws.Range("H36").Value = '47%'
$ws.Range("J36").Value = '1015.3 hPa'
$ws.Range("N36").Value = '11.3 °C 4:20 TU'
$ws.Range("O36").Value = '13.1 °C'
$ws.Range("E37").Value = '2026-02-17 04:49:54'
$ws.Range("H37").Value = '36%'
$ws.Range("J37").Value = '1015.0 hPa'
$ws.Range("N37").Value = '5.2 °C 4:27 TU'
$ws.Range("O37").Value = '8.6 °C'
$ws.Range("E38").Value = '2026-02-17 04:49:56'
$ws.Range("O38").Value = '8.9 °C'
$ws.Range("E39").Value = '2026-02-17 04:49:59'
$ws.Range("H39").Value = '57%'
$ws.Range("I39").Value = '0.9 mm'
$ws.Range("L39").Value = '108.0 km/h - 355º 4:12 TU'
$ws.Range("M39").Value = '-4.2 °C 4:29 TU'
$ws.Range("O39").Value = '-5.0 °C'
$ws.Range("E40").Value = '2026-02-17 04:50:01'
$ws.Range("J40").Value = '1017.3 hPa'
$ws.Range("O40").Value = '5.5 °C'
$ws.Range("E41").Value = '2026-02-17 04:50:03'
$ws.Range("H41").Value = '48%'
$ws.Range("K41").Value = '-0.1 MJ/m2'
$ws.Range("N41").Value = '13.8 °C 4:29 TU'
$ws.Range("O41").Value = '15.3 °C'
$ws.Range("E42").Value = '2026-02-17 04:50:06'
$ws.Range("E43").Value = '2026-02-17 04:50:08'
$ws.Range("N43").Value = '3.1 °C 4:19 TU'
$ws.Range("O43").Value = '4.6 °C'
$ws.Range("E44").Value = '2026-02-17 04:50:10'
$ws.Range("H44").Value = '67%'
$ws.Range("E45").Value = '2026-02-17 04:50:13'
$ws.Range("H45").Value = '52%'
$ws.Range("J45").Value = '1020.9 hPa'
$ws.Range("O45").Value = '4.5 °C'
$ws.Range("E46").Value = '2026-02-17 04:50:15'
$ws.Range("J46").Value = '1017.5 hPa'
$ws.Range("N46").Value = '12.6 °C 4:06 TU'
$ws.Range("O46").Value = '13.6 °C'
